$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 27 was blank; fill it in with a new "Git Ignore" entry (git log graph /
# git ignore system additions referenced in the commit message).
$ws.Range("A27").Value = "Git Ignore"
$ws.Range("B27").Value = "git ignore system"
$ws.Range("C27").Value = "1. Edit {home}/.git/info/exclude to ignore files, Every line regex the file that should (not) be ignored:`n# *.java    <- ignore javas globally`n# !*.java  <- don’t ignore java globally`n2. Put .gitignore file in any folder to state the ignore target in that folder.`nREF-- https://git-scm.com/docs/gitignore"

# The new row wraps onto several lines; match the taller row height.
$ws.Rows.Item(27).RowHeight = 64.5

# Move the active selection from B27 to C27.
$ws.Range("C27").Select()
